# Update the "registerSubscriptions" sheet: the bad-request test rows
# (rows 7-10) used to report a generic "Illegal Subscription sentence"
# rspMessage with rspCode 107001. They now report the more specific
# "Illegal GraphQL input" rspMessage with rspCode 101301.

$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("registerSubscriptions")

for ($r = 7; $r -le 10; $r++) {
    $wsRegister.Cells.Item($r, 7).Value = 101301
    $wsRegister.Cells.Item($r, 8).Value = "Illegal GraphQL input"
}

# Restore the previously-selected cells on each sheet.
$wsRegister.Range("D2").Select()

$wsDelete = $wb.Worksheets.Item("delSubscriptionById")
$wsDelete.Range("A5").Select()

$wsDelete.Activate()
